# Ensure correctness of tma2 examples
#
# - Rename sheet "PTEN" -> "pten"
# - Move the active/selected tab from the "pten" sheet to the "TMA map" sheet
# - Strip the (redundant) explicit style from a handful of cells on the
#   "pten" sheet so they fall back to the default cell style

$wb = $excel.ActiveWorkbook

$tmaMap = $wb.Worksheets.Item(1)
$pten   = $wb.Worksheets.Item(2)

# Sheet name "PTEN" -> "pten"
$pten.Name = "pten"

# These cells previously carried an explicit (but redundant) style index;
# put them back on the workbook's default ("Normal") style.
$cellsToReset = @(
    "C2", "J2", "K2",
    "J3", "K3",
    "E7", "F7", "I7", "J7", "K7",
    "C8", "D8", "F8", "K8", "L8"
)
foreach ($addr in $cellsToReset) {
    $pten.Range($addr).Style = "Normal"
}

# "TMA map" becomes the selected/active sheet (it previously was not),
# while "pten" is no longer the active tab.
$tmaMap.Activate()
